# no-op test
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
